# Update "想去人数" (interest count) values in column F across all four sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    "展览" = @{
        2 = 2274
        3 = 343
        4 = 180
        5 = 184
        6 = 345
        8 = 708
        9 = 516
        10 = 670
        11 = 372
        13 = 361
        15 = 4701
        16 = 161
        17 = 15
        18 = 24
        19 = 251
        20 = 143
        21 = 109
        23 = 90
        25 = 261
        26 = 106
    }
    "演出" = @{
        4 = 150
        6 = 182
        7 = 208
        8 = 2769
        14 = 107
        16 = 2505
    }
    "本地生活" = @{
        3 = 42
        4 = 402
        5 = 174
    }
    "全部类型" = @{
        3 = 42
        6 = 2274
        7 = 402
        8 = 343
        9 = 180
        10 = 184
        11 = 345
        12 = 150
        15 = 182
        16 = 174
        17 = 708
        18 = 516
        19 = 670
        20 = 372
        22 = 361
        24 = 4701
        25 = 208
        26 = 2769
        30 = 161
        31 = 15
        32 = 24
        35 = 251
        36 = 143
        37 = 109
        39 = 107
        41 = 90
        43 = 261
        44 = 106
        45 = 2505
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Cells.Item($row, 6).Value = $rows[$row]
    }
}
